$wb = $excel.ActiveWorkbook

# --- Sheets lookup ---
$wsModulos          = $wb.Worksheets.Item("03_MODULOS")
$wsFuncionalidades  = $wb.Worksheets.Item("04_MODULOS_FUNCIONALIDADES")

# --- Add the new trailing sheet "Hoja1" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHoja1 = $wb.Worksheets.Add($null, $lastSheet)
$wsHoja1.Name = "Hoja1"

# --- Populate new sheet header / data (order matters for shared-string ids) ---
$wsHoja1.Range("A1").Value = "CLASE"
$wsHoja1.Range("B1").Value = "NAMESPACE"
$wsFuncionalidades.Range("D4").Value = "PersonaNatural"
$wsHoja1.Range("A2").Value = "PersonaNatural"
$wsHoja1.Range("B2").Value = "JusNucleo.Bl.Personas"
$wsHoja1.Range("C1").Value = "FUCIONALIDAD_CODIGO"
$wsHoja1.Range("C2").Value = "FUNPERSONA"
$wsFuncionalidades.Range("E4").Value = "/Index"

# --- Column widths on 04_MODULOS_FUNCIONALIDADES (D/E got wider to fit new content) ---
$wsFuncionalidades.Columns.Item(4).ColumnWidth = 14.7109375
$wsFuncionalidades.Columns.Item(5).ColumnWidth = 9.42578125

# --- Selections / active cell bookkeeping ---
$wsModulos.Range("B2").Select()
$wsFuncionalidades.Range("E5").Select()
$wsHoja1.Range("C3").Select()

# --- Make 04_MODULOS_FUNCIONALIDADES the active/visible tab ---
$wsFuncionalidades.Activate()

Write-Host "done"
